$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The D column holds numeric-looking text (e.g. "30.532.65", "1.004") that must
# stay plain text (it was stored as inline text in the workbook). Pre-format the
# whole data range as Text before assigning values so Excel does not coerce them
# into real numbers, then restore the default "Normal" style so no new styling
# is introduced on the cells.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.532.65'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '1.920.08'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '239.34'
$ws.Range('E5').Value = '  -2.32%  '
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').Value = '0.4794'
$ws.Range('E7').Value = '  -1.67%  '
$ws.Range('D8').Value = '0.2882'
$ws.Range('E8').Value = '  -2.57%  '
$ws.Range('D9').Value = '0.06708'
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').Value = '18.79'
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('D11').Value = '104.31'
$ws.Range('E11').Value = '  -2.55%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.927.25'
$ws.Range('E12').Value = '  -1.50%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07741'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').Value = '5.264'
$ws.Range('E14').Value = '  -4.13%  '
$ws.Range('D15').Value = '0.6841'
$ws.Range('E15').Value = '  -2.74%  '
$ws.Range('D16').Value = '265.97'
$ws.Range('E16').Value = '  -6.47%  '
$ws.Range('D17').Value = '30.588.98'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '0.000007524'
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('D20').Value = '12.74'
$ws.Range('E20').Value = '  -3.44%  '
$ws.Range('D21').Value = '5.451'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('D22').Value = '1.004'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('B23').Value = 'BitDAO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D23').Value = '0.4567'
$ws.Range('E23').Value = '  -8.90%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '6.360'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = '9.676'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '163.75'
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '19.08'
$ws.Range('E27').Value = '  -4.49%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.091'
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '0.1025'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.388'
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.657'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '1.518'
$ws.Range('E32').Value = '  -4.21%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '4.271'
$ws.Range('E33').Value = '  -3.69%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04770'
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.7405'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.121'
$ws.Range('E36').Value = '  -4.38%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.685'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01948'
$ws.Range('E39').Value = '  -3.11%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.642'
$ws.Range('E40').Value = '  -2.16%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '6.351'
$ws.Range('E41').Value = '  -2.97%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '75.82'
$ws.Range('E42').Value = '  -1.78%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '2.009'
$ws.Range('E43').Value = '  -4.95%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.8625'
$ws.Range('E44').Value = '  -2.89%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '106.09'
$ws.Range('E45').Value = '  -2.54%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').Value = '0.4304'
$ws.Range('E46').Value = '  -3.67%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '7.593'
$ws.Range('E48').Value = '  -6.38%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '998.30'
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.1205'
$ws.Range('E50').Value = '  -4.23%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = '35.25'
$ws.Range('E51').Value = '  -1.85%  '

$ws.Range('D2:D51').Style = 'Normal'
